$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 32
$ws.Range("H32").Value = 450.33334
$ws.Range("I32").Value = 433.33334
$ws.Range("J32").Value = 467.33334
$ws.Range("K32").Value = 433.33334
$ws.Range("L32").Value = 467.33334
$ws.Range("M32").Value = -107.33334
$ws.Range("N32").Value = -1119.33334
# row 55
$ws.Range("H55").Value = 180.1
$ws.Range("I55").Value = 79.5
$ws.Range("J55").Value = 205.25
$ws.Range("K55").Value = 79.5
$ws.Range("L55").Value = 205.25
$ws.Range("M55").Value = 134.5
$ws.Range("N55").Value = -633.25
# row 64
$ws.Range("H64").Value = 2933.3333
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 2900
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 2900
$ws.Range("M64").Value = -2752
$ws.Range("N64").Value = -3396
# row 67
$ws.Range("H67").Value = 2933.3333
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 2900
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 2900
$ws.Range("M67").Value = -2142
$ws.Range("N67").Value = -4616
# row 98
$ws.Range("H98").Value = 6330.5
$ws.Range("I98").Value = 3383.1667
$ws.Range("J98").Value = 10751.5
$ws.Range("K98").Value = 3383.1667
$ws.Range("L98").Value = 10751.5
$ws.Range("M98").Value = -1885.1667
$ws.Range("N98").Value = -13747.5
# row 121
$ws.Range("H121").Value = 2246.52
$ws.Range("J121").Value = 2325.5417
$ws.Range("L121").Value = 6976.625100000001
$ws.Range("N121").Value = -10470.6251
# row 122
$ws.Range("H122").Value = 6330.5
$ws.Range("I122").Value = 3383.1667
$ws.Range("J122").Value = 10751.5
$ws.Range("K122").Value = 10149.5001
$ws.Range("L122").Value = 32254.5
$ws.Range("M122").Value = -7699.500100000001
$ws.Range("N122").Value = -37154.5
# row 137
$ws.Range("H137").Value = 2659.1738
$ws.Range("I137").Value = 1621.1945
$ws.Range("J137").Value = 6395.9
$ws.Range("K137").Value = 4863.583500000001
$ws.Range("L137").Value = 19187.7
$ws.Range("M137").Value = -2313.583500000001
$ws.Range("N137").Value = -24287.7

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 6
$ws.Range("H6").Value = 9871
$ws.Range("J6").Value = 10499.5
$ws.Range("L6").Value = 10499.5
$ws.Range("N6").Value = -10845.5
# row 32
$ws.Range("H32").Value = 10076.949
$ws.Range("I32").Value = 6326.492
$ws.Range("J32").Value = 16640.25
$ws.Range("K32").Value = 6326.492
$ws.Range("L32").Value = 16640.25
$ws.Range("M32").Value = -6039.492
$ws.Range("N32").Value = -17214.25
# row 61
$ws.Range("H61").Value = 1595.7273
$ws.Range("I61").Value = 1237.8462
$ws.Range("J61").Value = 2925
$ws.Range("K61").Value = 1237.8462
$ws.Range("L61").Value = 2925
$ws.Range("M61").Value = -1025.8462
$ws.Range("N61").Value = -3349
# row 74
$ws.Range("H74").Value = 1673.6285
$ws.Range("I74").Value = 1177.6
$ws.Range("J74").Value = 2913.7
$ws.Range("K74").Value = 1177.6
$ws.Range("L74").Value = 2913.7
$ws.Range("M74").Value = -303.5999999999999
$ws.Range("N74").Value = -4661.7
# row 77
$ws.Range("H77").Value = 1673.6285
$ws.Range("I77").Value = 1177.6
$ws.Range("J77").Value = 2913.7
$ws.Range("K77").Value = 5888
$ws.Range("L77").Value = 14568.5
$ws.Range("M77").Value = -1520
$ws.Range("N77").Value = -23304.5
# row 136
$ws.Range("H136").Value = 1595.7273
$ws.Range("I136").Value = 1237.8462
$ws.Range("J136").Value = 2925
$ws.Range("K136").Value = 3713.5386
$ws.Range("L136").Value = 8775
$ws.Range("M136").Value = -1163.5386
$ws.Range("N136").Value = -13875

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 94
$ws.Range("H94").Value = 1603.6666
$ws.Range("I94").Value = 1493.3334
$ws.Range("J94").Value = 1861.1111
$ws.Range("K94").Value = 1493.3334
$ws.Range("L94").Value = 1861.1111
$ws.Range("M94").Value = -1042.3334
$ws.Range("N94").Value = -2763.1111
# row 140
$ws.Range("H140").Value = 48306.25
$ws.Range("J140").Value = 48306.25
$ws.Range("L140").Value = 48306.25
$ws.Range("N140").Value = -58666.25
# row 141
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").ClearContents()
$ws.Range("N141").Value = 0

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Range("H16").Value = 3969411.2
$ws.Range("I16").Value = 6945474.5
$ws.Range("K16").Value = 6945474.5
$ws.Range("M16").Value = -6945187.5
# row 105
$ws.Range("H105").Value = 3456.7778
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 3638.875
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 3638.875
$ws.Range("M105").Value = -253
$ws.Range("N105").Value = -7132.875
# row 113
$ws.Range("H113").Value = 3969411.2
$ws.Range("I113").Value = 6945474.5
$ws.Range("K113").Value = 6945474.5
$ws.Range("M113").Value = -6943304.5
# row 132
$ws.Range("H132").Value = 2352.0625
$ws.Range("I132").Value = 2059.25
$ws.Range("J132").Value = 3230.5
$ws.Range("K132").Value = 6177.75
$ws.Range("L132").Value = 9691.5
$ws.Range("M132").Value = -3647.75
$ws.Range("N132").Value = -14751.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 64
$ws.Range("H64").Value = 9513.429
$ws.Range("J64").Value = 9513.429
$ws.Range("L64").Value = 28540.287
$ws.Range("N64").Value = -29080.287
# row 67
$ws.Range("H67").Value = 9513.429
$ws.Range("J67").Value = 9513.429
$ws.Range("L67").Value = 28540.287
$ws.Range("N67").Value = -30412.287
# row 92
$ws.Range("H92").Value = 26318326
$ws.Range("I92").Value = 525.8570999999999
$ws.Range("J92").Value = 41670376
$ws.Range("K92").Value = 1577.5713
$ws.Range("L92").Value = 125011128
$ws.Range("M92").Value = -329.5712999999998
$ws.Range("N92").Value = -125013624
# row 129
$ws.Range("H129").Value = 2378.4146
$ws.Range("J129").Value = 2375.4
$ws.Range("L129").Value = 7126.200000000001
$ws.Range("N129").Value = -17126.2
# row 132
$ws.Range("H132").Value = 3718.7715
$ws.Range("I132").Value = 1119.3334
$ws.Range("J132").Value = 5075
$ws.Range("K132").Value = 10074.0006
$ws.Range("L132").Value = 45675
$ws.Range("M132").Value = -7544.000599999999
$ws.Range("N132").Value = -50735
# row 136
$ws.Range("H136").Value = 3039.0908
$ws.Range("I136").Value = 2828.75
$ws.Range("J136").Value = 3600
$ws.Range("K136").Value = 8486.25
$ws.Range("L136").Value = 10800
$ws.Range("M136").Value = -3386.25
$ws.Range("N136").Value = -21000

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 57
$ws.Range("H57").Value = 33800
$ws.Range("J57").Value = 36666.668
$ws.Range("L57").Value = 36666.668
$ws.Range("N57").Value = -38306.668
# row 70
$ws.Range("H70").Value = 5646.8
$ws.Range("J70").Value = 5931.1875
$ws.Range("L70").Value = 5931.1875
$ws.Range("N70").Value = -6471.1875
# row 73
$ws.Range("H73").Value = 5646.8
$ws.Range("J73").Value = 5931.1875
$ws.Range("L73").Value = 5931.1875
$ws.Range("N73").Value = -7803.1875

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 16
$ws.Range("H16").Value = 578.8182
$ws.Range("I16").Value = 578.8182
$ws.Range("K16").Value = 578.8182
$ws.Range("M16").Value = -408.8182
# row 68
$ws.Range("H68").Value = 876.5
$ws.Range("I68").Value = 769.4828
$ws.Range("J68").Value = 3980
$ws.Range("K68").Value = 769.4828
$ws.Range("L68").Value = 3980
$ws.Range("M68").Value = -20.4828
$ws.Range("N68").Value = -5478
# row 71
$ws.Range("H71").Value = 876.5
$ws.Range("I71").Value = 769.4828
$ws.Range("J71").Value = 3980
$ws.Range("K71").Value = 3847.414
$ws.Range("L71").Value = 19900
$ws.Range("M71").Value = -103.4139999999998
$ws.Range("N71").Value = -27388

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 108
$ws.Range("H108").Value = 39800
$ws.Range("J108").Value = 39800
$ws.Range("L108").Value = 39800
$ws.Range("N108").Value = -47480
# row 122
$ws.Range("H122").Value = 4467.773
$ws.Range("I122").Value = 2962.818
$ws.Range("K122").Value = 8888.454000000002
$ws.Range("M122").Value = -6438.454000000002
# row 132
$ws.Range("H132").Value = 7095969
$ws.Range("I132").Value = 4246.9287
$ws.Range("K132").Value = 12740.7861
$ws.Range("M132").Value = -10210.7861
# row 135
$ws.Range("H135").Value = 76553.766
$ws.Range("J135").Value = 76553.766
$ws.Range("L135").Value = 76553.766
$ws.Range("N135").Value = -86693.766
